# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 322
$ws1.Range("F3").Value = 78
$ws1.Range("F4").Value = 484
$ws1.Range("F5").Value = 4719
$ws1.Range("F6").Value = 372
$ws1.Range("F8").Value = 290
$ws1.Range("F9").Value = 731
$ws1.Range("F10").Value = 211

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 322
$ws4.Range("F3").Value = 78
$ws4.Range("F4").Value = 484
$ws4.Range("F5").Value = 4719
$ws4.Range("F6").Value = 372
$ws4.Range("F8").Value = 290
$ws4.Range("F9").Value = 731
$ws4.Range("F11").Value = 211
